# Commit: "Did extra planning and prepared for gameplay code
#          started game concept art"
#
# Adds a planning paragraph that explains the reasoning behind turning the
# pause menu into its own internal window/surface instead of refreshing the
# whole screen. It goes right after the "Own Surface" Heading2, into the
# paragraph that already carries the <w:noProof/> mark (an empty paragraph
# left there for a screenshot/placeholder).

$d = $word.ActiveDocument

$p = $d.Paragraphs.Item(16)
$r = $p.Range
$r.Select()

$sel = $word.Selection
$sel.Collapse(0)

$newText = "Initially I made the menus refresh the entire screen. This was easy to implement but meant that the game would not be visible around the pause menu. I wanted the menus to be their own internal window that I could place anywhere on the screen. "

$sel.TypeText($newText)

# Re-select exactly what was just typed and mark it as "no proofing" so the
# run formatting matches the rest of this (image-caption style) paragraph.
$sel.MoveStart(1, -$newText.Length)
$sel.NoProofing = $true
